# Scoreboard.xlsx edit script
# - Swap the Minute2("G") values between ScoreM rows 2 and 5 (team 1 <-> team 4)
# - Fill in new Workout 6 columns (L, M, N) on the ScoreF sheet for rows 2-10
# - Move the active sheet from ScoreM to ScoreF, and update each sheet's
#   remembered selection (ScoreM -> G6, ScoreF -> B25, FM -> C16)

$wb = $excel.ActiveWorkbook

$wsScoreM = $wb.Worksheets.Item("ScoreM")
$wsScoreF = $wb.Worksheets.Item("ScoreF")
$wsFM     = $wb.Worksheets.Item("FM")

# --- ScoreM: swap G2 and G5 ---
$wsScoreM.Range("G2").Value = 35.1
$wsScoreM.Range("G5").Value = 35

# --- ScoreF: add Workout 6 points/time columns (L:N) for rows 2-10 ---
$scoreFData = @{
    2  = @(37, 38, 10000)
    3  = @(43, 13, 10000)
    4  = @(41, 1,  10000)
    5  = @(38, 4,  10000)
    6  = @(39, 45, 10000)
    7  = @(40, 30, 10000)
    8  = @(38, 10, 10000)
    9  = @(38, 47, 10000)
    10 = @(37, 42, 10000)
}

foreach ($row in $scoreFData.Keys) {
    $vals = $scoreFData[$row]
    $wsScoreF.Range("L$row").Value = $vals[0]
    $wsScoreF.Range("M$row").Value = $vals[1]
    $wsScoreF.Range("N$row").Value = $vals[2]
}

# --- Update remembered selection per sheet, and set ScoreF as the active tab ---
$wsScoreM.Activate()
$wsScoreM.Range("G6").Select()

$wsFM.Activate()
$wsFM.Range("C16").Select()

$wsScoreF.Activate()
$wsScoreF.Range("B25").Select()
